# quote puller done, now working on shuffle
#
# Wrap every quote in column A (rows 2-51 of the "quotes" sheet) in single
# quote characters, e.g. `Don't explain your philosophy. Embody it.`
# becomes `'Don't explain your philosophy. Embody it.'`.
#
# A literal leading apostrophe can't be written through `.Value` directly
# (Excel treats a leading `'` as the text-prefix marker rather than literal
# content), so each cell is staged via `.Formula` (a proper quoted string
# literal) and then flattened back down to a plain value with
# Copy/PasteSpecial(xlPasteValues) so the saved cell stays a normal shared
# string instead of a formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 51 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $orig = $cell.Value2
    if ($orig -eq $null) { continue }

    $quoted = "'" + $orig + "'"
    $escaped = $quoted.Replace('"', '""')

    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$excel.CutCopyMode = 0

# Column B (author) ends up auto-sized to its content once the sheet is
# touched again.
$ws.Columns("B").ColumnWidth = 13.43

# Reset the lingering multi-column selection left over from the prior
# editing session back to the top-left cell.
$ws.Range("A1").Select() | Out-Null
